$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell AO1: text "38", same style as AN1 (bold/centered) ---
$ws.Range("AN1").Copy($ws.Range("AO1"))
$ws.Range("BZ1").Formula = '="38"'
$ws.Range("BZ1").Copy()
$ws.Range("AO1").PasteSpecial(-4163)
$ws.Range("BZ1").Clear()

# --- Week 38 (column AO) data values ---
$ws.Range("AO2").Value = 0
$ws.Range("AO5").Value = 0
$ws.Range("AO6").Value = 1
$ws.Range("AO7").Value = 0
$ws.Range("AO8").Value = 0
$ws.Range("AO9").Value = 0
$ws.Range("AO10").Value = 0
$ws.Range("AO11").Value = 0
$ws.Range("AO14").Value = 0
$ws.Range("AO15").Value = 0
$ws.Range("AO16").Value = 0
$ws.Range("AO17").Value = 0
$ws.Range("AO18").Value = 0
$ws.Range("AO19").Value = 0
$ws.Range("AO23").Value = 0
$ws.Range("AO25").Value = 0
$ws.Range("AO28").Value = 0
$ws.Range("AO29").Value = 1
$ws.Range("AO31").Value = 0
$ws.Range("AO35").Value = 6
$ws.Range("AO36").Value = 0
$ws.Range("AO37").Value = 0
$ws.Range("AO38").Value = 0
$ws.Range("AO41").Value = 0
$ws.Range("AO42").Value = 0
$ws.Range("AO43").Value = 0
$ws.Range("AO44").Value = 0
$ws.Range("AO45").Value = 0
$ws.Range("AO46").Value = 0
$ws.Range("AO47").Value = 0
$ws.Range("AO48").Value = 0
$ws.Range("AO49").Value = 0
$ws.Range("AO50").Value = 0
$ws.Range("AO51").Value = 0
$ws.Range("AO52").Value = 0
$ws.Range("AO53").Value = 0
$ws.Range("AO54").Value = 0
$ws.Range("AO55").Value = 0
$ws.Range("AO56").Value = 0
$ws.Range("AO57").Value = 0
$ws.Range("AO58").Value = 0

$excel.CutCopyMode = $false
